# Updated IND model - 2025-08-07 07:10
#
# Target sheet: VEDA_Sets-Proc (the ~TFM_Psets table).
#  - Row 3  (CCGT set):     add ",*GasCC*" to the PSET_PN pattern (B3),
#                           duplicate the SetName into SetDesc (G3),
#                           and add the T_Pos_AndOr / T_Neg_AndOr flags (H3/I3).
#  - Row 7  (OCGT set):     add ",EN*CT*" to the PSET_PN pattern (B7),
#                           and add the T_Pos_AndOr / T_Neg_AndOr flags (H7/I7).
#  - Row 17 (Nuclear set):  add a new PSET_PN exclusion pattern "-*SMR" (B17),
#                           and add the T_Pos_AndOr / T_Neg_AndOr flags (H17/I17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Row 3 - CCGT
$ws.Range("B3").Value = "ep_gas_combined_cycle*,ep_oil_combined_cycle*,CCGT*,*GasCC*"
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("H3").Value = "And"
$ws.Range("I3").Value = "Or"

# Row 7 - OCGT (Peaker)
$ws.Range("B7").Value = "ep_gas_gas_turbine*,ep_oil_gas_turbine*,gas turbine*,EN*CT*"
$ws.Range("H7").Value = "And"
$ws.Range("I7").Value = "Or"

# Row 17 - Nuclear
$ws.Range("B17").Value = "-*SMR"
$ws.Range("H17").Value = "And"
$ws.Range("I17").Value = "Or"
